$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J column values
$ws.Range("J3").Value = 1.1399999999999999
$ws.Range("J5").Value = 1.23
$ws.Range("J8").Value = 1.04
$ws.Range("J9").Value = 0.85
$ws.Range("J15").Value = 1.0900000000000001

# Update K column formulas (value recalculated automatically)
$ws.Range("K3").Formula = "=0.48/2"
$ws.Range("K8").Formula = "=0.61/2"
$ws.Range("K9").Formula = "=0.49/2"
$ws.Range("K15").Formula = "=0.93/2"

# Update selected cell in the sheet view
$ws.Range("K16").Select()
